# Adds a new "2022-Q3" sheet (fund holdings detail) positioned right after
# the "总计" summary sheet, pushing the existing quarter sheets down, and
# updates the "总计" sheet with the new quarter's roll-up row.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: shift the quarterly figures down
#    one row and insert the new 2022-Q3 totals at the top of the table.
# ---------------------------------------------------------------------
$summary = $sheets.Item(1)

# Push the existing rows 2-7 down to rows 3-8 (content only - column A's
# running index is left untouched since it already lines up row-for-row).
# NB: read with .Value2 - .Value's getter is unreliable in this host when
# its result is consumed directly (by Write-Host or as an assignment RHS).
$summary.Range("B8").Value2 = $summary.Range("B7").Value2
$summary.Range("C8").Value2 = $summary.Range("C7").Value2
$summary.Range("D8").Value2 = $summary.Range("D7").Value2

$summary.Range("B7").Value2 = $summary.Range("B6").Value2
$summary.Range("C7").Value2 = $summary.Range("C6").Value2
$summary.Range("D7").Value2 = $summary.Range("D6").Value2

$summary.Range("B6").Value2 = $summary.Range("B5").Value2
$summary.Range("C6").Value2 = $summary.Range("C5").Value2
$summary.Range("D6").Value2 = $summary.Range("D5").Value2

$summary.Range("B5").Value2 = $summary.Range("B4").Value2
$summary.Range("C5").Value2 = $summary.Range("C4").Value2
$summary.Range("D5").Value2 = $summary.Range("D4").Value2

$summary.Range("B4").Value2 = $summary.Range("B3").Value2
$summary.Range("C4").Value2 = $summary.Range("C3").Value2
$summary.Range("D4").Value2 = $summary.Range("D3").Value2

$summary.Range("B3").Value2 = $summary.Range("B2").Value2
$summary.Range("C3").Value2 = $summary.Range("C2").Value2
$summary.Range("D3").Value2 = $summary.Range("D2").Value2

# New row-8 index cell (copy style from the row above, then set its value).
$summary.Range("A7").Copy($summary.Range("A8"))
$summary.Range("A8").Value = 6

# New 2022-Q3 totals go into row 2.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 13
$summary.Range("D2").Value = 4.58

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" detail sheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $sheets.Add($null, $summary)
$q3.Name = "2022-Q3"

function Set-HeaderCell($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Bold = $true
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
}

function Set-IndexCell($ws, $addr, $num) {
    $c = $ws.Range($addr)
    $c.Value = $num
    $c.Font.Bold = $true
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
}

function Set-TextCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
}

function Set-NumberCell($ws, $row, $col, $num) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $num
}

# Header row.
Set-HeaderCell $q3 "B1" "基金代码"
Set-HeaderCell $q3 "C1" "基金名称"
Set-HeaderCell $q3 "D1" "基金规模"
Set-HeaderCell $q3 "E1" "股票总仓位"
Set-HeaderCell $q3 "F1" "仓位占比"
Set-HeaderCell $q3 "G1" "持有市值(亿元)"
Set-HeaderCell $q3 "H1" "仓位排名"

# Fund rows: idx, code, name, size, position, ratio, marketValue, rank.
# marketValue is written as text EXCEPT for the two zero-size funds, where
# the source data stored a literal numeric 0.
$rows = @(
    @(0,  "001606", "农银汇理工业4.0灵活配置混合",   "40.20", "86.85", "4.34", "1.7447", 8),
    @(1,  "000336", "农银研究精选混合",               "38.76", "82.60", "2.86", "1.1085", 9),
    @(2,  "001182", "易方达安心回馈混合A",             "52.92", "35.21", "1.12", "0.5927", 9),
    @(3,  "213008", "宝盈资源优选混合",                "9.37",  "85.08", "5.24", "0.4910", 6),
    @(4,  "213006", "宝盈核心优势灵活配置混合A",        "9.16",  "73.82", "3.99", "0.3655", 7),
    @(5,  "010383", "宝盈基础产业混合A",               "3.00",  "89.54", "3.95", "0.1185", 9),
    @(6,  "000259", "农银区间收益混合",                "4.20",  "69.93", "1.95", "0.0819", 6),
    @(7,  "001707", "诺安高端制造股票A",               "1.12",  "90.57", "4.09", "0.0458", 4),
    @(8,  "005638", "农银汇理量化智慧动力混合",         "0.60",  "88.21", "2.79", "0.0167", 3),
    @(9,  "010384", "宝盈基础产业混合C",               "0.22",  "89.54", "3.95", "0.0087", 9),
    @(10, "000241", "宝盈核心优势灵活配置混合C",        "0.21",  "73.82", "3.99", "0.0084", 7),
    @(11, "014536", "诺安高端制造股票C",               "0.00",  "90.57", "4.09", 0,        4),
    @(12, "016594", "易方达安心回馈混合C",              "0.00",  "35.21", "1.12", 0,        9)
)

$r = 2
foreach ($row in $rows) {
    $addr = "A" + $r
    Set-IndexCell $q3 $addr $row[0]
    Set-TextCell  $q3 $r 2 $row[1]
    Set-TextCell  $q3 $r 3 $row[2]
    Set-TextCell  $q3 $r 4 $row[3]
    Set-TextCell  $q3 $r 5 $row[4]
    Set-TextCell  $q3 $r 6 $row[5]
    $gVal = $row[6]
    if ($gVal -is [string]) {
        Set-TextCell $q3 $r 7 $gVal
    } else {
        Set-NumberCell $q3 $r 7 $gVal
    }
    Set-NumberCell $q3 $r 8 $row[7]
    $r++
}
